$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row at position 13, shifting rows 13-21 down to 14-22.
$ws.Rows("13:13").Insert()

# 2) Update row 10 (Objetivos: body) with the new objectives text.
$objectivesText = "Desenvolver um projeto sobre tema de Engenharia de Produção, similar a situações que os alunos irão encontrar na vida real no efetivo exercício de sua profissão, `nAplicar e integrar conhecimentos adquiridos em demais disciplinas de seu curso`nDesenvolver competências técnicas, as relacionadas ao projeto em si, bem como competências transversais (habilidades e atitudes), num ambiente de aprendizagem baseado em PBL (Project-Baed Learning)."
$ws.Range("B10").Value = $objectivesText
$ws.Range("C10").Value = $objectivesText

# 3) New row 13 holds the "Docentes responsáveis" value (no label in column A).
$ws.Range("B10:C10").Copy()
$ws.Range("B13:C13").PasteSpecial(-4122)
$teacherText = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Range("B13").Value = $teacherText
$ws.Range("C13").Value = $teacherText
$ws.Range("A13").Clear()

# 4) Row 14 (Programa resumido:) body text.
$shortSyllabusText = "Tópicos que abordem o tema do projeto de seu planejamento a execução."
$ws.Range("B14").Value = $shortSyllabusText
$ws.Range("C14").Value = $shortSyllabusText

# 5) Row 16 (Programa:) body text.
$programText = "Noções de Gestão de Projetos`nOrganização do tempo: dimensão pessoal;`nTécnicas para a realização de apresentações;`nNoções de Aprendizagem Baseada em Projetos`nTrabalho em Grupo, Equipes e times. `nPostura e Ética Profissional`nTécnicas para redação de relatório técnico;`nTutoria de projetos.`nAssuntos Técnicos específicos relacionados com o tema do projeto."
$ws.Range("B16").Value = $programText
$ws.Range("C16").Value = $programText

# 6) Row 19 (Método:) body text.
$methodText = "O método utilizado tem por fundamento a Aprendizagem Baseada em Projetos (PBL) que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, aspectos de liderança e capacidade de comunicação, dentre outras.`n`nOs alunos serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a um tema de Engenharia de Produção, similar ao que eles irão encontrar na vida real no efetivo exercício de sua profissão. `nCada grupo deverá buscar o conhecimento prático necessário para ser aplicado no desenvolvimento do projeto.`nAs aulas ocorrerão: 1) através de uma reunião da equipe de trabalho para tratar do projeto, e  2) palestras e dinâmicas relativas ao tema do projeto, conduzidas por professores  ou profissionais de empresas."
$ws.Range("B19").Value = $methodText
$ws.Range("C19").Value = $methodText

# 7) Row 20 (Critério:) body text.
$criteriaText = "A nota será individual e será a média ponderada de componentes do projeto, tais como: Projeto Preliminar, Projeto Final, envolvimento do aluno com o projeto, Avaliação dos Pares, Apresentação de Trabalhos, dentre outros.`nO detalhamento dos pesos para ponderação da média da disciplina será definido por uma equipe de professores que atuarão na coordenação da disciplina."
$ws.Range("B20").Value = $criteriaText
$ws.Range("C20").Value = $criteriaText

# 8) Row 21 (Norma de recuperação:) body text.
$recoveryText = "Não há recuperação"
$ws.Range("B21").Value = $recoveryText
$ws.Range("C21").Value = $recoveryText

# 9) New row 22 (Bibliografia: + body text).
$ws.Range("A21:C21").Copy()
$ws.Range("A22:C22").PasteSpecial(-4122)
$biblioLabel = "Bibliografia:"
$biblioText = "Artigos sobre metodologias ativas de aprendizagem e  Project Based Learning.`nLivros e Artigos científicos relacionados com o tema do projeto."
$ws.Range("A22").Value = $biblioLabel
$ws.Range("B22").Value = $biblioText
$ws.Range("C22").Value = $biblioText

$ws.Range("A1").Select()
